$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 3628.080185983

$ws.Range("B3").Value = -868.5473439574538
$ws.Range("C3").Value = 0.4813056583245627
$ws.Range("D3").Value = 3624.397065511

$ws.Range("C4").Value = 0.40476687146591733
$ws.Range("D4").Value = 3622.680581468

$ws.Range("B5").Value = -888.6086061064475
$ws.Range("C5").Value = 0.6716507717197328
$ws.Range("D5").Value = 3626.814940312

$ws.Range("D6").Value = 3633.756014323

$ws.Range("C7").Value = 0.6159822546319933
$ws.Range("D7").Value = 3628.025892173

$ws.Range("C8").Value = 0.35733090501693315
$ws.Range("D8").Value = 3625.271143391

$ws.Range("C9").Value = 0.22204873711219358
$ws.Range("D9").Value = 3624.784063875

$ws.Range("B10").Value = -857.5898287715161
$ws.Range("C10").Value = 1.1829585622111276
$ws.Range("D10").Value = 3621.851209648

$ws.Range("D11").Value = 3623.151370983
